$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sfb")

# Row 8 - Digital - (UV)
$ws.Range("H8").Value = 179317
$ws.Range("J8").Value = 4657659
$ws.Range("K8").Value = 320973
$ws.Range("M8").Value = 8265551
$ws.Range("N8").Value = 234395
$ws.Range("O8").Value = 29.73
$ws.Range("P8").Value = 6967932

# Row 11 - ZZ Total
$ws.Range("H11").Value = 312827
$ws.Range("J11").Value = 6754605
$ws.Range("K11").Value = 600446
$ws.Range("M11").Value = 12520842
$ws.Range("N11").Value = 234395
$ws.Range("O11").Value = 29.73
$ws.Range("P11").Value = 6967932
